$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose "Value" column (D) was empty and now gets "DNI"
$dniRows = @(18, 19, 20, 21, 22, 23, 24, 41)

foreach ($r in $dniRows) {
    $ws.Range("D$r").Value = "DNI"
}

# Reflect the cell selection left active by the author when the file was saved
$ws.Range("D12").Select()
